$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Value corrections in the H:M (2nd/summary) columns across the three task blocks
$ws.Range("I3").Value = 2
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 5
$ws.Range("M3").Value = 3

$ws.Range("I6").Value = 5

$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 2
$ws.Range("L10").Value = 12
$ws.Range("M10").Value = 8

$ws.Range("I14").Value = 5
$ws.Range("L14").Value = 17

$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 2
$ws.Range("L18").Value = 12
$ws.Range("M18").Value = 8

$ws.Range("J20").Value = 5
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 4
$ws.Range("M20").Value = 7

$ws.Range("I22").Value = 5
$ws.Range("L22").Value = 11
